# Trade #51 closed at 2026-02-17 21:08:22 - unknown UNKNOWN +0.000%
#
# Updates the live-trading workbook to reflect:
#  1. A recalculated Summary roll-up (capital / pnl / trade counts / win rate).
#  2. A recalculated MarketMaking row in the Strategy Status roll-up.
#  3. Trade #79 (row 80 on "All Trades", row 47 on "MarketMaking") being
#     closed out (CLOSED, exit price, pnl, capital-after, exit reason,
#     duration).
#  4. A brand new open trade, #112, appended to both the "All Trades" and
#     "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.61
$summary.Range("B4").Value = 0.41
$summary.Range("B6").Value = 79
$summary.Range("B7").Value = 38
$summary.Range("B9").Value = 48.1

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.61
$status.Range("D5").Value = 46
$status.Range("E5").Value = 0.3
$status.Range("F5").Value = 0.61
$status.Range("G5").Value = 52.17

# ---------------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out trade #79 (row 80): Exit Price, Status, P&L %, P&L $,
# Capital After, Exit Reason, Duration (min)
$allTrades.Range("G80").Value = 0.98
$allTrades.Range("H80").Value = "CLOSED"
$allTrades.Range("I80").Value = 1.0309
$allTrades.Range("J80").Value = 0.01
$allTrades.Range("K80").Value = 100.61
$allTrades.Range("L80").Value = "early_exit"
$allTrades.Range("M80").Value = 0.14

# Append new trade #112 as row 113
$allTrades.Range("A113").Value = 112
$allTrades.Range("B113").NumberFormat = "@"
$allTrades.Range("B113").Value = "2026-02-17"
$allTrades.Range("C113").Value = "21:08:12"
$allTrades.Range("D113").Value = "MarketMaking"
$allTrades.Range("E113").Value = "UP"
$allTrades.Range("F113").Value = 0.97
$allTrades.Range("G113").Formula = "=""" + """"
$allTrades.Range("H113").Value = "OPEN"
$allTrades.Range("I113").Value = 0
$allTrades.Range("J113").Value = 0
$allTrades.Range("K113").Value = 100.6014872031006
$allTrades.Range("L113").Formula = "=""" + """"
$allTrades.Range("M113").Value = 0
$allTrades.Range("N113").Value = 0
$allTrades.Range("O113").Value = 0
$allTrades.Range("P113").Value = 0.6
$allTrades.Range("Q113").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out trade #79 (row 47): Exit Price, Status, P&L %, P&L $,
# Capital After, Exit Reason, Duration (min)
$mm.Range("G47").Value = 0.98
$mm.Range("H47").Value = "CLOSED"
$mm.Range("I47").Value = 1.0309
$mm.Range("J47").Value = 0.01
$mm.Range("K47").Value = 100.61
$mm.Range("P47").Value = "early_exit"
$mm.Range("Q47").Value = 0.14

# Append new trade #112 as row 80
$mm.Range("A80").Value = 112
$mm.Range("B80").NumberFormat = "@"
$mm.Range("B80").Value = "2026-02-17"
$mm.Range("C80").Value = "21:08:12"
$mm.Range("D80").Value = "MarketMaking"
$mm.Range("E80").Value = "UP"
$mm.Range("F80").Value = 0.97
$mm.Range("G80").Formula = "=""" + """"
$mm.Range("H80").Value = "OPEN"
$mm.Range("I80").Value = 0
$mm.Range("J80").Value = 0
$mm.Range("K80").Value = 100.6014872031006
$mm.Range("L80").Value = 0
$mm.Range("M80").Value = 0
$mm.Range("N80").Value = 0.6
$mm.Range("O80").Value = "Normal spread capture: 19600 bps"
$mm.Range("P80").Formula = "=""" + """"
$mm.Range("Q80").Value = 0
